$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New numeric values (inherit default/no explicit style) ---
$ws.Range("J4").Value = 32

$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 3

$ws.Range("G9").Value = 1

$ws.Range("C11").Value = 32
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 5
$ws.Range("H11").Value = 3

# --- New blank cells that carry the grey "s=6" border/fill style ---
# Copy the format from an existing s="6" cell onto each new blank cell.
$ws.Range("I5").Copy()
$ws.Range("J5").PasteSpecial(-4122)

$ws.Range("I6").Copy()
$ws.Range("J6").PasteSpecial(-4122)

$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)

$ws.Range("I8").Copy()
$ws.Range("J8").PasteSpecial(-4122)

$ws.Range("I9").Copy()
$ws.Range("J9").PasteSpecial(-4122)

$ws.Range("I10").Copy()
$ws.Range("J10").PasteSpecial(-4122)

$ws.Range("I5").Copy()
$ws.Range("G11").PasteSpecial(-4122)

$ws.Range("I5").Copy()
$ws.Range("J11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Update the active selection to match the new commit state ---
$ws.Range("I17").Select()
